# Update ResumoInscricoes figures for 2025-1 Superior intake.
# Each touched row gets its "Inscritos" (E) bumped by the newly-homologated
# registration(s); where pagamento/isenção totals also shifted, "Pagos" (F)
# and "Inscrições homologadas" (H) are bumped to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 22

$ws.Range("E15").Value = 145

$ws.Range("E19").Value = 46
$ws.Range("F19").Value = 24
$ws.Range("H19").Value = 24

$ws.Range("E32").Value = 18

$ws.Range("F38").Value = 14
$ws.Range("H38").Value = 14

$ws.Range("E42").Value = 30

$ws.Range("E43").Value = 21

$ws.Range("E46").Value = 24

$ws.Range("E48").Value = 26
$ws.Range("F48").Value = 15
$ws.Range("H48").Value = 15

$ws.Range("E61").Value = 26

$ws.Range("E64").Value = 31

$ws.Range("E70").Value = 38
$ws.Range("F70").Value = 17
$ws.Range("H70").Value = 17

$ws.Range("E73").Value = 28
$ws.Range("F73").Value = 11
$ws.Range("H73").Value = 11

$ws.Range("E78").Value = 41
$ws.Range("F78").Value = 17
$ws.Range("H78").Value = 17

$ws.Range("E79").Value = 30
